# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" sheets represents the source file
# fecc3372-687d-4d6d-801d-3d3a89bc3f86.md. A handback was produced for it,
# but it turned out to not be the latest version, so:
#   - "Latest Target File" (I7) and "Latest Handback File" (J7) get filled in
#   - "Latest Handback DateTime" (K7) gets a real timestamp
#   - "Error Detail" (P7) gets a "not the latest" warning message
#   - I7 becomes a hyperlink (like the already-populated A7 cell)

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/69e09abb1877563d6752fbe1e916444865d35693/e2e/fecc3372-687d-4d6d-801d-3d3a89bc3f86.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/646c5592be94d7852a86a9c38f8cee0346aa8b61/e2e/fecc3372-687d-4d6d-801d-3d3a89bc3f86.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("I7").Value = "fecc3372-687d-4d6d-801d-3d3a89bc3f86.md"
$wsZhCn.Range("I7").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/69e09abb1877563d6752fbe1e916444865d35693/e2e/fecc3372-687d-4d6d-801d-3d3a89bc3f86.md", "", "", "fecc3372-687d-4d6d-801d-3d3a89bc3f86.md")

$wsZhCn.Range("J7").Value = "fecc3372-687d-4d6d-801d-3d3a89bc3f86.6b321e5955bd71f29db922e13f1d629c50a4be4c.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-17 16:54:27"
$wsZhCn.Range("P7").Value = $errorMessage

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("I7").Value = "fecc3372-687d-4d6d-801d-3d3a89bc3f86.md"
$wsDeDe.Range("I7").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/69e09abb1877563d6752fbe1e916444865d35693/e2e/fecc3372-687d-4d6d-801d-3d3a89bc3f86.md", "", "", "fecc3372-687d-4d6d-801d-3d3a89bc3f86.md")

$wsDeDe.Range("J7").Value = "fecc3372-687d-4d6d-801d-3d3a89bc3f86.6b321e5955bd71f29db922e13f1d629c50a4be4c.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-17 16:54:35"
$wsDeDe.Range("P7").Value = $errorMessage
